$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.635.49'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = '  -1.90%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.588.87'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = '  -2.37%  '

$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("E6").Value = '  -3.21%  '

$ws.Range("E7").Value = '  +0.14%  '

$ws.Range("E8").Value = '  -2.62%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0614'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = '  -2.14%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.57'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '  -4.15%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.810.89'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = '  -2.39%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.587.64'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  -2.60%  '

$ws.Range("E14").Value = '  -3.16%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.521'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '  -4.50%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.69'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  -0.22%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.611.29'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  -1.95%  '

$ws.Range("E18").Value = '  -2.53%  '

$ws.Range("E19").Value = '  +0.27%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '207.94'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '  -4.30%  '

$ws.Range("E21").Value = '  -3.49%  '

$ws.Range("E22").Value = '  -3.50%  '

$ws.Range("E23").Value = '  -3.96%  '

$ws.Range("E24").Value = '  -2.45%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.78'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = '  -0.71%  '

$ws.Range("E26").Value = '  +0.11%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.23'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = '  -1.25%  '

$ws.Range("E28").Value = '  -3.93%  '

$ws.Range("E29").Value = '  -2.39%  '

$ws.Range("E30").Value = '  -0.46%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.14'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = '  -2.34%  '

$ws.Range("E32").Value = '  -4.38%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.661'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = '  +20.00%  '

$ws.Range("E34").Value = '  -3.05%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.305.36'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = '  -3.26%  '

$ws.Range("E36").Value = '  -1.26%  '

$ws.Range("E37").Value = '  -5.59%  '

$ws.Range("E38").Value = '  -3.55%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.827'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = '  -3.44%  '

$ws.Range("E40").Value = '  +0.20%  '

$ws.Range("E41").Value = '  -1.73%  '

$ws.Range("E42").Value = '  +1.92%  '

$ws.Range("E43").Value = '  -3.38%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '62.58'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '  -4.52%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.724.15'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  -2.25%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '89.48'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = '  -1.28%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.61'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = '  -1.05%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.838'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = '  -1.73%  '

$ws.Range("E49").Value = '  -1.92%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0976'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = '  -2.35%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.47'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = '  -1.52%  '
